# Amend corrected label annotations
# Updates column F ("labels") values: lowercases label text and, for multi-label
# cells joined with ' || ', reverses the order of the individual label segments
# (matching the exact new values captured from the source diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(9, 6).Value = '93_referral_statement'
$ws.Cells.Item(15, 6).Value = '18_hazards_to_humans_and_domestic_animals'
$ws.Cells.Item(25, 6).Value = 'ppe'
$ws.Cells.Item(29, 6).Value = 'application instructions'
$ws.Cells.Item(30, 6).Value = 'env warning - water || off target movement'
$ws.Cells.Item(31, 6).Value = '32_physical_and_chemical_hazards'
$ws.Cells.Item(34, 6).Value = 'application instructions'
$ws.Cells.Item(35, 6).Value = 'application instructions'
$ws.Cells.Item(36, 6).Value = 'application instructions'
$ws.Cells.Item(37, 6).Value = '135_product_information'
$ws.Cells.Item(38, 6).Value = 'mixing || 135_product_information'
$ws.Cells.Item(45, 6).Value = 'use restrictions'
$ws.Cells.Item(46, 6).Value = 'use restrictions'
$ws.Cells.Item(48, 6).Value = 'off target movement'
$ws.Cells.Item(49, 6).Value = 'off target movement'
$ws.Cells.Item(50, 6).Value = 'off target movement'
$ws.Cells.Item(51, 6).Value = 'off target movement'
$ws.Cells.Item(52, 6).Value = 'application instructions || off target movement'
$ws.Cells.Item(53, 6).Value = 'off target movement'
$ws.Cells.Item(54, 6).Value = 'off target movement'
$ws.Cells.Item(55, 6).Value = 'off target movement'
$ws.Cells.Item(57, 6).Value = 'off target movement'
$ws.Cells.Item(59, 6).Value = 'off target movement'
$ws.Cells.Item(60, 6).Value = 'off target movement'
$ws.Cells.Item(61, 6).Value = 'off target movement'
$ws.Cells.Item(62, 6).Value = 'application instructions || off target movement'
$ws.Cells.Item(63, 6).Value = 'application instructions || off target movement'
$ws.Cells.Item(65, 6).Value = '172_sensitive_areas || application instructions || off target movement'
$ws.Cells.Item(66, 6).Value = 'application instructions'
$ws.Cells.Item(67, 6).Value = 'application instructions'
$ws.Cells.Item(68, 6).Value = 'application instructions'
$ws.Cells.Item(69, 6).Value = 'application instructions'
$ws.Cells.Item(71, 6).Value = 'application instructions'
$ws.Cells.Item(73, 6).Value = 'application instructions'
$ws.Cells.Item(74, 6).Value = 'application instructions'
$ws.Cells.Item(75, 6).Value = 'application instructions'
$ws.Cells.Item(78, 6).Value = 'mixing'
$ws.Cells.Item(79, 6).Value = 'mixing'
$ws.Cells.Item(80, 6).Value = 'mixing'
$ws.Cells.Item(81, 6).Value = 'mixing'
$ws.Cells.Item(82, 6).Value = 'mixing'
$ws.Cells.Item(84, 6).Value = 'mixing'
$ws.Cells.Item(85, 6).Value = 'mixing'
$ws.Cells.Item(86, 6).Value = 'mixing'
$ws.Cells.Item(88, 6).Value = 'safety procedures'
$ws.Cells.Item(90, 6).Value = 'safety procedures'
$ws.Cells.Item(92, 6).Value = 'mixing'
$ws.Cells.Item(95, 6).Value = 'application instructions'
$ws.Cells.Item(155, 6).Value = 'use restrictions'
$ws.Cells.Item(156, 6).Value = 'use restrictions'
$ws.Cells.Item(160, 6).Value = 'use restrictions'
$ws.Cells.Item(162, 6).Value = 'use restrictions'
$ws.Cells.Item(165, 6).Value = 'use restrictions'
$ws.Cells.Item(166, 6).Value = 'application instructions'
$ws.Cells.Item(168, 6).Value = 'use restrictions'
$ws.Cells.Item(169, 6).Value = 'application instructions'
$ws.Cells.Item(172, 6).Value = 'mixing'
$ws.Cells.Item(173, 6).Value = 'mixing'
$ws.Cells.Item(174, 6).Value = 'mixing'
$ws.Cells.Item(175, 6).Value = 'mixing'
$ws.Cells.Item(178, 6).Value = 'mixing'
$ws.Cells.Item(179, 6).Value = 'mixing'
$ws.Cells.Item(193, 6).Value = 'mixing'
$ws.Cells.Item(194, 6).Value = 'application instructions'
$ws.Cells.Item(196, 6).Value = 'mixing'
$ws.Cells.Item(197, 6).Value = 'use restrictions'
$ws.Cells.Item(198, 6).Value = 'application instructions'
$ws.Cells.Item(200, 6).Value = 'use restrictions'
$ws.Cells.Item(202, 6).Value = 'use restrictions'
$ws.Cells.Item(203, 6).Value = 'application instructions'
$ws.Cells.Item(204, 6).Value = 'application instructions'
$ws.Cells.Item(205, 6).Value = 'use restrictions'
$ws.Cells.Item(206, 6).Value = 'application instructions'
$ws.Cells.Item(208, 6).Value = 'application instructions'
$ws.Cells.Item(209, 6).Value = 'use restrictions'
$ws.Cells.Item(210, 6).Value = 'use restrictions'
$ws.Cells.Item(211, 6).Value = 'mixing'
$ws.Cells.Item(215, 6).Value = 'application instructions'
$ws.Cells.Item(223, 6).Value = 'mixing'
$ws.Cells.Item(226, 6).Value = 'application instructions'
$ws.Cells.Item(227, 6).Value = 'application instructions'
$ws.Cells.Item(228, 6).Value = 'use restrictions'
$ws.Cells.Item(229, 6).Value = 'irrigation || chemigation'
$ws.Cells.Item(230, 6).Value = 'safety procedures || chemigation'
$ws.Cells.Item(232, 6).Value = 'irrigation'
$ws.Cells.Item(234, 6).Value = 'application instructions || chemigation'
$ws.Cells.Item(235, 6).Value = 'safety procedures'
$ws.Cells.Item(236, 6).Value = 'use restrictions || irrigation'
$ws.Cells.Item(237, 6).Value = 'application instructions'
$ws.Cells.Item(238, 6).Value = 'use restrictions'
$ws.Cells.Item(242, 6).Value = 'use restrictions'
$ws.Cells.Item(243, 6).Value = 'application instructions'
$ws.Cells.Item(246, 6).Value = 'use restrictions'
$ws.Cells.Item(247, 6).Value = 'application instructions'
$ws.Cells.Item(249, 6).Value = 'application instructions'
$ws.Cells.Item(251, 6).Value = 'application instructions'
$ws.Cells.Item(253, 6).Value = 'application instructions'
$ws.Cells.Item(254, 6).Value = 'use restrictions'
$ws.Cells.Item(256, 6).Value = 'use restrictions'
$ws.Cells.Item(259, 6).Value = 'application instructions'
$ws.Cells.Item(260, 6).Value = 'use restrictions'
$ws.Cells.Item(261, 6).Value = 'application instructions'
$ws.Cells.Item(262, 6).Value = 'application instructions'
$ws.Cells.Item(263, 6).Value = 'use restrictions'
$ws.Cells.Item(265, 6).Value = 'use restrictions'
$ws.Cells.Item(267, 6).Value = '154_pesticide_storage'
$ws.Cells.Item(284, 6).Value = '93_referral_statement'
$ws.Cells.Item(290, 6).Value = '18_hazards_to_humans_and_domestic_animals'
$ws.Cells.Item(291, 6).Value = 'application instructions'
$ws.Cells.Item(292, 6).Value = 'env warning - water || off target movement'
$ws.Cells.Item(294, 6).Value = '154_pesticide_storage'
